# fix: reorder tab after add/remove contact (#234)
# Update the "Community" sheet's sample row (row 2) with:
#   - a new default membership setting in column B
#   - unit prices for meal/drink tickets in column D
#   - a refreshed updatedAt timestamp in column H
#   - the updatedBy contact email in column I

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Community")

$ws.Range("B2").Value = '{"membershipFee":null,"membershipEmail":null}'
$ws.Range("D2").Value = '[{"name":"meal","count":null,"unitPrice":"5.00","hidden":false},{"name":"drink","count":null,"unitPrice":"1.00","hidden":false},{"name":"cotton-candy","count":null,"unitPrice":null,"hidden":false}]'
$ws.Range("H2").Value = "2025-09-29T18:29:07.690Z"
$ws.Range("I2").Value = "dev@email.com"
